$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous last row (row 64) changes from the "latest" date-only format
# to the standard datetime format, since it's no longer the most recent entry.
$ws.Range("A64").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new daily row (row 65)
$ws.Range("A65").Value = 45652
$ws.Range("A65").NumberFormat = "YYYY-MM-DD"
$ws.Range("B65").Value = 155
$ws.Range("C65").Value = 144
$ws.Range("D65").Value = 149
